$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.227.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.050.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.22%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.539"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.047.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.83%  "
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("E14").Value = "  -4.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.119"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.558.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.322.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("E18").Value = "  -2.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.050.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "471.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.703"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.96%  "
$ws.Range("E23").Value = "  -1.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.74%  "
$ws.Range("E29").Value = "  +1.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("E31").Value = "  -2.34%  "
$ws.Range("E32").Value = "  -2.87%  "
$ws.Range("E33").Value = "  -1.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0820"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.48%  "
$ws.Range("E36").Value = "  -1.95%  "
$ws.Range("E37").Value = "  +1.83%  "
$ws.Range("E38").Value = "  -3.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "438.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.51%  "
$ws.Range("E43").Value = "  -2.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.94%  "
$ws.Range("E45").Value = "  +3.07%  "
$ws.Range("E46").Value = "  -4.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.784.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.53%  "
